$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.112.97"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "2.904.89"
$ws.Range("E3").Value = "  +8.00%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "196.77"
$ws.Range("E5").Value = "  +4.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "601.38"
$ws.Range("E6").Value = "  +1.95%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.194"
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("D10").Value = "2.904.06"
$ws.Range("E10").Value = "  +8.08%  "
$ws.Range("E11").Value = "  +10.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.161"
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.95"
$ws.Range("E13").Value = "  +4.04%  "
$ws.Range("D14").Value = "3.435.11"
$ws.Range("E14").Value = "  +7.83%  "
$ws.Range("D15").Value = "76.081.98"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.66"
$ws.Range("E17").Value = "  +3.85%  "
$ws.Range("D18").Value = "2.902.81"
$ws.Range("E18").Value = "  +7.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.01"
$ws.Range("E19").Value = "  -4.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.65"
$ws.Range("E20").Value = "  +4.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "384.38"
$ws.Range("E21").Value = "  +2.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.32"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.15"
$ws.Range("E24").Value = "  +2.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.29"
$ws.Range("E26").Value = "  +2.22%  "
$ws.Range("E27").Value = "  +7.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.85"
$ws.Range("E28").Value = "  +4.05%  "
$ws.Range("E29").Value = "  +14.71%  "
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "515.19"
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.85"
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("E34").Value = "  +3.62%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.06"
$ws.Range("E36").Value = "  +1.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.26"
$ws.Range("E37").Value = "  +4.74%  "
$ws.Range("E38").Value = "  -2.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.68"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "184.55"
$ws.Range("E40").Value = "  +7.93%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("E42").Value = "  +4.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.09"
$ws.Range("E43").Value = "  +0.99%  "
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0934"
$ws.Range("E45").Value = "  +10.25%  "
$ws.Range("E46").Value = "  +3.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.23"
$ws.Range("E47").Value = "  +2.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.40"
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("E49").Value = "  +8.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.680"
$ws.Range("E50").Value = "  +14.67%  "
$ws.Range("E51").Value = "  +3.06%  "
